$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = -1.160000000000096
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = -0.2100000000000648
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = -4.67000000000003
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = 48.24999999999996
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = 10.95999999999999
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = -10.23999999999998
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = -12.59
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = 3.339999999999989
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = 13.05999999999997
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = 71.71999999999998
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = 18.33000000000004
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = -6.119999999999976
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = 2.809999999999988
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = 17.92000000000006
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = -6.989999999999981
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = -2.530000000000044
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = 12.07999999999998
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = -5.519999999999982
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = 64.27999999999993
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = 13.09999999999989
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = -7.060000000000059
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = -7.810000000000031
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = 10.17000000000007
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = 30.81999999999996
